# Update cryptocurrency price/volume data refreshed by GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal numbers (e.g. "1.00", "40.70").
# Mark those cells as Text first so Excel keeps the exact string (with
# trailing zeros / precision) instead of silently re-typing them as numbers.
$numericPriceCells = @(
    'D4',
    'D5',
    'D6',
    'D7',
    'D10',
    'D11',
    'D12',
    'D15',
    'D17',
    'D19',
    'D20',
    'D22',
    'D23',
    'D24',
    'D26',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D41',
    'D47',
    'D48',
    'D50',
    'D51'
)
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.231.08'
$ws.Range('E2').Value = '  +1.97%  '

$ws.Range('D3').Value = '2.563.49'
$ws.Range('E3').Value = '  +1.48%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '318.04'
$ws.Range('E5').Value = '  +0.95%  '

$ws.Range('D6').Value = '97.38'
$ws.Range('E6').Value = '  +3.45%  '

$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  +0.58%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('E9').Value = '  +2.45%  '

$ws.Range('D10').Value = '35.84'
$ws.Range('E10').Value = '  +0.79%  '

$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  +1.25%  '

$ws.Range('D12').Value = '7.53'

$ws.Range('E13').Value = '  -4.79%  '

$ws.Range('D14').Value = '2.957.03'
$ws.Range('E14').Value = '  +1.50%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '15.13'
$ws.Range('E15').Value = '  -1.94%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.481.96'
$ws.Range('E16').Value = '  -1.72%  '

$ws.Range('D17').Value = '0.851'
$ws.Range('E17').Value = '  +0.89%  '

$ws.Range('D18').Value = '43.172.73'
$ws.Range('E18').Value = '  +1.63%  '

$ws.Range('D19').Value = '6.85'
$ws.Range('E19').Value = '  +4.48%  '

$ws.Range('D20').Value = '12.67'
$ws.Range('E20').Value = '  -1.35%  '

$ws.Range('D21').Value = '0.0₃0967'
$ws.Range('E21').Value = '  +0.97%  '

$ws.Range('D22').Value = '70.05'
$ws.Range('E22').Value = '  -0.73%  '

$ws.Range('D23').Value = '254.09'
$ws.Range('E23').Value = '  +1.85%  '

$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +1.31%  '

$ws.Range('E25').Value = '  +2.60%  '

$ws.Range('D26').Value = '27.02'
$ws.Range('E26').Value = '  +1.61%  '

$ws.Range('E27').Value = '  +0.33%  '

$ws.Range('E28').Value = '  +1.67%  '

$ws.Range('D29').Value = '40.70'
$ws.Range('E29').Value = '  +4.31%  '

$ws.Range('D30').Value = '10.31'
$ws.Range('E30').Value = '  +2.01%  '

$ws.Range('D31').Value = '5.89'
$ws.Range('E31').Value = '  -0.45%  '

$ws.Range('D32').Value = '156.27'
$ws.Range('E32').Value = '  +0.30%  '

$ws.Range('D33').Value = '19.31'
$ws.Range('E33').Value = '  -0.15%  '

$ws.Range('D34').Value = '2.13'
$ws.Range('E34').Value = '  +0.69%  '

$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '2.71'
$ws.Range('E35').Value = '  +3.48%  '

$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '3.34'
$ws.Range('E36').Value = '  +1.21%  '

$ws.Range('D37').Value = '0.0800'
$ws.Range('E37').Value = '  +2.50%  '

$ws.Range('D38').Value = '0.112'
$ws.Range('E38').Value = '  +1.77%  '

$ws.Range('D39').Value = '2.45'
$ws.Range('E39').Value = '  +4.71%  '

$ws.Range('E40').Value = '  +0.23%  '

$ws.Range('D41').Value = '22.31'
$ws.Range('E41').Value = '  -5.66%  '

$ws.Range('E42').Value = '  +1.16%  '

$ws.Range('E43').Value = '  +2.06%  '

$ws.Range('E45').Value = '  -1.51%  '

$ws.Range('D46').Value = '1.998.35'
$ws.Range('E46').Value = '  -1.04%  '

$ws.Range('D47').Value = '85.45'
$ws.Range('E47').Value = '  +1.49%  '

$ws.Range('D48').Value = '9.11'
$ws.Range('E48').Value = '  +3.43%  '

$ws.Range('D49').Value = '2.808.79'
$ws.Range('E49').Value = '  +1.50%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '74.95'
$ws.Range('E50').Value = '  +3.22%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '105.04'
$ws.Range('E51').Value = '  +2.98%  '
